$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 32, shifting existing rows 32:151 down to 33:152
$ws.Rows("32:32").Insert()

# Fill the new row 32 with its data (constant columns copied from the pattern,
# plus the new date/volume/price/origin values from the commit)
$ws.Cells.Item(32, 1).Value = 5
$ws.Cells.Item(32, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(32, 3).Value = "Maule"
$ws.Cells.Item(32, 4).Value = (Get-Date -Year 2023 -Month 11 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Cells.Item(32, 5).Value = 7
$ws.Cells.Item(32, 6).Value = 100112026
$ws.Cells.Item(32, 7).Value = "Haba"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 500
$ws.Cells.Item(32, 11).Value = 9000
$ws.Cells.Item(32, 12).Value = 9000
$ws.Cells.Item(32, 13).Value = 9000
$ws.Cells.Item(32, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(32, 15).Value = "Región del Maule"
$ws.Cells.Item(32, 16).Value = 360
$ws.Cells.Item(32, 17).Value = 25
$ws.Cells.Item(32, 18).Value = "Hortaliza"
